$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131 (shifts existing rows 131-201 down to 132-202)
$ws.Rows(131).Insert()

# Populate the new row 131 with the new weekly price record
$ws.Range("A131").Value = 1
$ws.Range("B131").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C131").Value = "Arica y Parinacota"
$ws.Range("D131").Value = 45141
$ws.Range("E131").Value = 15
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100108
$ws.Range("H131").Value = "Tropicales y subtropicales"
$ws.Range("I131").Value = 100108003
$ws.Range("J131").Value = "Maracuyá"
$ws.Range("K131").Value = "Sin especificar"
$ws.Range("L131").Value = "Primera"
$ws.Range("M131").Value = 130
$ws.Range("N131").Value = 20000
$ws.Range("O131").Value = 22000
$ws.Range("P131").Value = 21000
$ws.Range("Q131").Value = "`$/caja 20 kilos"
$ws.Range("R131").Value = "Región de Arica y Parinacota"
$ws.Range("S131").Value = 1050
$ws.Range("T131").Value = 20
